# "Added last minute updates"
#
# The first paragraph of the document (the hidden **ID__...__ID** marker
# paragraph) gets:
#   1. A paragraph border added (pBdr) with 5-twip spacing on all 4 sides,
#      matching the style already used by the other paragraphs below it.
#   2. Its left indent bumped from 120 -> 225 twips (6pt -> 11.25pt),
#      again matching the other paragraphs.
#   3. Its marker text updated from
#        **ID__AFFARS_5351_topic_1__ID**
#      to
#        **ID__AFFARS_PART_5351__ID**
#      and the trailing run containing a lone space is dropped, leaving a
#      single run with just the marker text.

$d = $word.ActiveDocument
$para = $d.Paragraphs(1)

# --- paragraph formatting -------------------------------------------------
# Add the paragraph border (<w:pBdr><w:top w:space="5"/> ... ) seen on the
# other paragraphs in this document.
$borders = $para.Range.ParagraphFormat.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromRight = 5

# Left indent 120 twips (6pt) -> 225 twips (11.25pt)
$para.Range.ParagraphFormat.LeftIndent = 11.25

# --- run/text content -------------------------------------------------
# Replace the paragraph's text (everything except the trailing paragraph
# mark) with the new marker text. This both updates the id string and
# removes the separate trailing " " run, leaving one run with the new text.
$r = $para.Range
[void]$r.MoveEnd(1, -1)
$r.Text = "**ID__AFFARS_PART_5351__ID**"
